# Commit: "Fruta / hortaliza, semanal"
# A new weekly observation row is inserted at row 108 of the single sheet,
# pushing the existing rows 108-216 down to 109-217 (dimension grows from
# A1:T216 to A1:T217). Populate the newly inserted row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 108; everything below shifts down by one.
$ws.Rows(108).Insert()

# Fill in the new row 108 with the new observation's data.
$ws.Cells.Item(108, 1).Value = 11
$ws.Cells.Item(108, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(108, 3).Value = "Bíobío"
$ws.Cells.Item(108, 4).Value = 44880
$ws.Cells.Item(108, 5).Value = 8
$ws.Cells.Item(108, 6).Value = "Fruta"
$ws.Cells.Item(108, 7).Value = 100108
$ws.Cells.Item(108, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(108, 9).Value = 100108005
$ws.Cells.Item(108, 10).Value = "Piña"
$ws.Cells.Item(108, 11).Value = "Caramelo"
$ws.Cells.Item(108, 12).Value = "Segunda"
$ws.Cells.Item(108, 13).Value = 350
$ws.Cells.Item(108, 14).Value = 24000
$ws.Cells.Item(108, 15).Value = 25000
$ws.Cells.Item(108, 16).Value = 24429
$ws.Cells.Item(108, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(108, 18).Value = "Ecuador"
$ws.Cells.Item(108, 19).Value = 1745
$ws.Cells.Item(108, 20).Value = 14
